# Tugas_1/Kelas.xlsx edit script
# 1. Menambahkan random function  -> incremental "No" formulas in column A (rows 3-11)
# 2. Menambahkan data agar spreadsheet tidak monoton -> vary Alamat/Pekerjaan/Alasan per row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A small pool of "profile" data (Alamat, Pekerjaan, Alasan) used to make the
# previously-monotonous D:F columns vary from row to row.
$alamat = @(
  "Jl. Kebon Jeruk No. 123",
  "Jl. Merdeka No. 45",
  "Jl. Cihampelas No. 67",
  "Jl. Pahlawan No. 89",
  "Jl. Ahmad Yani No. 12",
  "Jl. Sukajadi No. 34",
  "Jl. Riau No. 56",
  "Jl. Dago No. 78",
  "Jl. Setiabudi No. 90",
  "Jl. Cijerah No. 23"
)
$pekerjaan = @(
  "Software Engineer",
  "Data Analyst",
  "UI/UX Designer",
  "Network Engineer",
  "System Administrator",
  "Software Developer",
  "IT Consultant",
  "Cloud Engineer",
  "DevOps Engineer",
  "Cyber Security Analyst"
)
$alasan = @(
  "Menyukai tantangan baru",
  "Ingin mengembangkan skill programming",
  "Minat pada desain grafis",
  "Passionate tentang jaringan komputer",
  "Pengalaman dalam administrasi sistem",
  "Minat pada pengembangan aplikasi",
  "Memiliki keahlian konsultasi IT",
  "Berpengalaman dalam pengelolaan cloud",
  "Memiliki keterampilan DevOps",
  "Minat pada keamanan cyber"
)

# Fill D:F for every data row (2-51), cycling through the 10 profiles so the
# sheet no longer repeats the same "Fresh Graduate Akademi Kominfo" / "Scalable
# Web Service with Golang (Batch1)" / "Kesempatan belajar golang" on every row.
for ($r = 2; $r -le 51; $r++) {
    $i = ($r - 2) % 10
    $ws.Range("D$r").Value = $alamat[$i]
    $ws.Range("E$r").Value = $pekerjaan[$i]
    $ws.Range("F$r").Value = $alasan[$i]
}

# Add an incrementing formula down column A for a handful of rows (A3:A11),
# continuing the numbering sequence (=row above + 1) instead of the previous
# hard-coded numbers.
$ws.Range("A3").Formula = "=A2+1"
$ws.Range("A4:A11").Formula = "=A3+1"

# Move the active selection/viewport to C9 (was D44, viewport anchored at A37).
$ws.Range("A1").Select()
$ws.Range("C9").Select()
